$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = 44372
$ws.Cells.Item(2, 13).Value = 60
$ws.Cells.Item(2, 14).Value = 20000
$ws.Cells.Item(2, 15).Value = 21000
$ws.Cells.Item(2, 16).Value = 20667
$ws.Cells.Item(2, 19).Value = 1033

# Row 3
$ws.Cells.Item(3, 4).Value = 44830
$ws.Cells.Item(3, 13).Value = 200
$ws.Cells.Item(3, 14).Value = 30000
$ws.Cells.Item(3, 15).Value = 30000
$ws.Cells.Item(3, 16).Value = 30000
$ws.Cells.Item(3, 19).Value = 1500

# Row 4
$ws.Cells.Item(4, 4).Value = 44365
$ws.Cells.Item(4, 13).Value = 150
$ws.Cells.Item(4, 14).Value = 20000
$ws.Cells.Item(4, 16).Value = 20500
$ws.Cells.Item(4, 19).Value = 1025

# Row 5
$ws.Cells.Item(5, 4).Value = 44357
$ws.Cells.Item(5, 13).Value = 200
$ws.Cells.Item(5, 14).Value = 20000
$ws.Cells.Item(5, 15).Value = 21000
$ws.Cells.Item(5, 16).Value = 20500
$ws.Cells.Item(5, 19).Value = 1025

# Row 6
$ws.Cells.Item(6, 4).Value = 44410
$ws.Cells.Item(6, 13).Value = 40
$ws.Cells.Item(6, 14).Value = 25000
$ws.Cells.Item(6, 15).Value = 25000
$ws.Cells.Item(6, 16).Value = 25000
$ws.Cells.Item(6, 19).Value = 1250

# Row 7
$ws.Cells.Item(7, 4).Value = 44355
$ws.Cells.Item(7, 13).Value = 200
$ws.Cells.Item(7, 14).Value = 20000
$ws.Cells.Item(7, 15).Value = 21000
$ws.Cells.Item(7, 16).Value = 20500
$ws.Cells.Item(7, 18).Value = "Ecuador"
$ws.Cells.Item(7, 19).Value = 1025

# Row 9
$ws.Cells.Item(9, 4).Value = 44277
$ws.Cells.Item(9, 13).Value = 60
$ws.Cells.Item(9, 14).Value = 24000
$ws.Cells.Item(9, 15).Value = 24000
$ws.Cells.Item(9, 16).Value = 24000
$ws.Cells.Item(9, 19).Value = 1200

# Row 10
$ws.Cells.Item(10, 4).Value = 44356
$ws.Cells.Item(10, 13).Value = 100
$ws.Cells.Item(10, 14).Value = 20000
$ws.Cells.Item(10, 15).Value = 21000
$ws.Cells.Item(10, 16).Value = 20500
$ws.Cells.Item(10, 19).Value = 1025

# Row 11
$ws.Cells.Item(11, 4).Value = 44442
$ws.Cells.Item(11, 13).Value = 30
$ws.Cells.Item(11, 14).Value = 22000
$ws.Cells.Item(11, 15).Value = 22000
$ws.Cells.Item(11, 16).Value = 22000
$ws.Cells.Item(11, 19).Value = 1100

# Row 12
$ws.Cells.Item(12, 4).Value = 44613
$ws.Cells.Item(12, 14).Value = 30000
$ws.Cells.Item(12, 15).Value = 30000
$ws.Cells.Item(12, 16).Value = 30000
$ws.Cells.Item(12, 19).Value = 1500

# Row 13
$ws.Cells.Item(13, 4).Value = 44424
$ws.Cells.Item(13, 13).Value = 70
$ws.Cells.Item(13, 14).Value = 24000
$ws.Cells.Item(13, 15).Value = 25000
$ws.Cells.Item(13, 16).Value = 24429
$ws.Cells.Item(13, 19).Value = 1221

# Row 14
$ws.Cells.Item(14, 4).Value = 44473
$ws.Cells.Item(14, 13).Value = 40

# Row 15
$ws.Cells.Item(15, 4).Value = 44298
$ws.Cells.Item(15, 13).Value = 240
$ws.Cells.Item(15, 14).Value = 19000
$ws.Cells.Item(15, 15).Value = 20000
$ws.Cells.Item(15, 16).Value = 19500
$ws.Cells.Item(15, 19).Value = 975

# Row 16
$ws.Cells.Item(16, 4).Value = 44300
$ws.Cells.Item(16, 13).Value = 150
$ws.Cells.Item(16, 16).Value = 19500
$ws.Cells.Item(16, 19).Value = 975

# Row 17
$ws.Cells.Item(17, 4).Value = 44326
$ws.Cells.Item(17, 13).Value = 40
$ws.Cells.Item(17, 14).Value = 22000
$ws.Cells.Item(17, 15).Value = 22000
$ws.Cells.Item(17, 16).Value = 22000
$ws.Cells.Item(17, 19).Value = 1100

# Row 18
$ws.Cells.Item(18, 4).Value = 44529
$ws.Cells.Item(18, 13).Value = 34
$ws.Cells.Item(18, 14).Value = 28000
$ws.Cells.Item(18, 15).Value = 28000
$ws.Cells.Item(18, 16).Value = 28000
$ws.Cells.Item(18, 19).Value = 1400

# Row 19
$ws.Cells.Item(19, 4).Value = 44435
$ws.Cells.Item(19, 13).Value = 60
$ws.Cells.Item(19, 14).Value = 25000
$ws.Cells.Item(19, 15).Value = 25000
$ws.Cells.Item(19, 16).Value = 25000
$ws.Cells.Item(19, 19).Value = 1250

# Row 20
$ws.Cells.Item(20, 4).Value = 44438
$ws.Cells.Item(20, 13).Value = 25
$ws.Cells.Item(20, 14).Value = 21000
$ws.Cells.Item(20, 15).Value = 21000
$ws.Cells.Item(20, 16).Value = 21000
$ws.Cells.Item(20, 19).Value = 1050

# Row 21
$ws.Cells.Item(21, 4).Value = 44166
$ws.Cells.Item(21, 13).Value = 120
$ws.Cells.Item(21, 14).Value = 28000
$ws.Cells.Item(21, 15).Value = 28000
$ws.Cells.Item(21, 16).Value = 28000
$ws.Cells.Item(21, 19).Value = 1400

# Row 22
$ws.Cells.Item(22, 4).Value = 44302
$ws.Cells.Item(22, 13).Value = 100
$ws.Cells.Item(22, 14).Value = 19000
$ws.Cells.Item(22, 15).Value = 20000
$ws.Cells.Item(22, 16).Value = 19500
$ws.Cells.Item(22, 19).Value = 975

# Row 23
$ws.Cells.Item(23, 4).Value = 44165
$ws.Cells.Item(23, 13).Value = 300
$ws.Cells.Item(23, 14).Value = 27000
$ws.Cells.Item(23, 16).Value = 27500
$ws.Cells.Item(23, 19).Value = 1375

# Row 24
$ws.Cells.Item(24, 4).Value = 44270
$ws.Cells.Item(24, 13).Value = 50
$ws.Cells.Item(24, 14).Value = 24000
$ws.Cells.Item(24, 15).Value = 24000
$ws.Cells.Item(24, 16).Value = 24000
$ws.Cells.Item(24, 19).Value = 1200

# Row 25
$ws.Cells.Item(25, 4).Value = 44363
$ws.Cells.Item(25, 13).Value = 150
$ws.Cells.Item(25, 14).Value = 21000
$ws.Cells.Item(25, 15).Value = 22000
$ws.Cells.Item(25, 16).Value = 21500
$ws.Cells.Item(25, 19).Value = 1075

# Row 26
$ws.Cells.Item(26, 4).Value = 44284
$ws.Cells.Item(26, 13).Value = 40
$ws.Cells.Item(26, 14).Value = 23000
$ws.Cells.Item(26, 15).Value = 23000
$ws.Cells.Item(26, 16).Value = 23000
$ws.Cells.Item(26, 19).Value = 1150

# Row 27
$ws.Cells.Item(27, 4).Value = 44350
$ws.Cells.Item(27, 13).Value = 90
$ws.Cells.Item(27, 14).Value = 21000
$ws.Cells.Item(27, 15).Value = 22000
$ws.Cells.Item(27, 16).Value = 21556
$ws.Cells.Item(27, 19).Value = 1078

# Row 28
$ws.Cells.Item(28, 4).Value = 44620
$ws.Cells.Item(28, 13).Value = 60
$ws.Cells.Item(28, 14).Value = 22000
$ws.Cells.Item(28, 15).Value = 22000
$ws.Cells.Item(28, 16).Value = 22000
$ws.Cells.Item(28, 19).Value = 1100

# Row 29
$ws.Cells.Item(29, 4).Value = 44382
$ws.Cells.Item(29, 13).Value = 200
$ws.Cells.Item(29, 16).Value = 19500
$ws.Cells.Item(29, 19).Value = 975

# Row 30
$ws.Cells.Item(30, 4).Value = 44305
$ws.Cells.Item(30, 13).Value = 40
$ws.Cells.Item(30, 14).Value = 24000
$ws.Cells.Item(30, 15).Value = 24000
$ws.Cells.Item(30, 16).Value = 24000
$ws.Cells.Item(30, 19).Value = 1200

# Row 31
$ws.Cells.Item(31, 4).Value = 44452
$ws.Cells.Item(31, 13).Value = 35
$ws.Cells.Item(31, 14).Value = 21000
$ws.Cells.Item(31, 15).Value = 22000
$ws.Cells.Item(31, 16).Value = 21429
$ws.Cells.Item(31, 19).Value = 1071

# Row 32
$ws.Cells.Item(32, 4).Value = 44445
$ws.Cells.Item(32, 13).Value = 35
$ws.Cells.Item(32, 14).Value = 20000
$ws.Cells.Item(32, 15).Value = 20000
$ws.Cells.Item(32, 16).Value = 20000
$ws.Cells.Item(32, 19).Value = 1000

# Row 33
$ws.Cells.Item(33, 4).Value = 44417
$ws.Cells.Item(33, 13).Value = 30
$ws.Cells.Item(33, 14).Value = 24000
$ws.Cells.Item(33, 15).Value = 24000
$ws.Cells.Item(33, 16).Value = 24000
$ws.Cells.Item(33, 19).Value = 1200

# Row 34
$ws.Cells.Item(34, 4).Value = 44354
$ws.Cells.Item(34, 13).Value = 150
$ws.Cells.Item(34, 14).Value = 21000
$ws.Cells.Item(34, 15).Value = 22000
$ws.Cells.Item(34, 16).Value = 21500
$ws.Cells.Item(34, 19).Value = 1075

# Row 35
$ws.Cells.Item(35, 4).Value = 44396
$ws.Cells.Item(35, 13).Value = 45
$ws.Cells.Item(35, 14).Value = 22000
$ws.Cells.Item(35, 15).Value = 22000
$ws.Cells.Item(35, 16).Value = 22000
$ws.Cells.Item(35, 19).Value = 1100

# Row 36
$ws.Cells.Item(36, 4).Value = 44760
$ws.Cells.Item(36, 13).Value = 300
$ws.Cells.Item(36, 14).Value = 24000
$ws.Cells.Item(36, 16).Value = 24500
$ws.Cells.Item(36, 19).Value = 1225

# Row 37
$ws.Cells.Item(37, 4).Value = 44522
$ws.Cells.Item(37, 12).Value = "Primera"
$ws.Cells.Item(37, 13).Value = 25
$ws.Cells.Item(37, 14).Value = 30000
$ws.Cells.Item(37, 15).Value = 30000
$ws.Cells.Item(37, 16).Value = 30000
$ws.Cells.Item(37, 19).Value = 1500

# Row 38
$ws.Cells.Item(38, 4).Value = 44431
$ws.Cells.Item(38, 13).Value = 60
$ws.Cells.Item(38, 14).Value = 25000
$ws.Cells.Item(38, 15).Value = 25000
$ws.Cells.Item(38, 16).Value = 25000
$ws.Cells.Item(38, 18).Value = "Perú"
$ws.Cells.Item(38, 19).Value = 1250

# Row 39
$ws.Cells.Item(39, 4).Value = 44312
$ws.Cells.Item(39, 13).Value = 50
$ws.Cells.Item(39, 14).Value = 22000
$ws.Cells.Item(39, 15).Value = 22000
$ws.Cells.Item(39, 16).Value = 22000
$ws.Cells.Item(39, 19).Value = 1100

# Row 40
$ws.Cells.Item(40, 4).Value = 44263
$ws.Cells.Item(40, 12).Value = "Segunda"
$ws.Cells.Item(40, 13).Value = 150
$ws.Cells.Item(40, 14).Value = 15000
$ws.Cells.Item(40, 15).Value = 15000
$ws.Cells.Item(40, 16).Value = 15000
$ws.Cells.Item(40, 19).Value = 750

# Row 41
$ws.Cells.Item(41, 4).Value = 44333
$ws.Cells.Item(41, 13).Value = 30
$ws.Cells.Item(41, 14).Value = 22000
$ws.Cells.Item(41, 16).Value = 22000
$ws.Cells.Item(41, 19).Value = 1100
